# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before "总计"),
#    populated with the fund-holding detail rows for that quarter.
# 2) Prepend a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item(2)          # "2021-Q4" - used as insertion anchor + style donor
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Header row - copy the header formatting (style index) from the "2021-Q4" sheet,
# then overwrite the text in each header cell.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Column-A index cells (0, 1, ...) reuse the bold/centered style used elsewhere.
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A3").PasteSpecial(-4122)

# Row 2: 005901 诺安汇利灵活配置混合A
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'005901"
$newSheet.Cells.Item(2,2).Style = "Normal"
$newSheet.Cells.Item(2,3).Value = "诺安汇利灵活配置混合A"
$newSheet.Cells.Item(2,4).Value = "'0.08"
$newSheet.Cells.Item(2,4).Style = "Normal"
$newSheet.Cells.Item(2,5).Value = "'86.88"
$newSheet.Cells.Item(2,5).Style = "Normal"
$newSheet.Cells.Item(2,6).Value = "'4.99"
$newSheet.Cells.Item(2,6).Style = "Normal"
$newSheet.Cells.Item(2,7).Value = "'0.0040"
$newSheet.Cells.Item(2,7).Style = "Normal"
$newSheet.Cells.Item(2,8).Value = 9

# Row 3: 005902 诺安汇利灵活配置混合C
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'005902"
$newSheet.Cells.Item(3,2).Style = "Normal"
$newSheet.Cells.Item(3,3).Value = "诺安汇利灵活配置混合C"
$newSheet.Cells.Item(3,4).Value = "'0.02"
$newSheet.Cells.Item(3,4).Style = "Normal"
$newSheet.Cells.Item(3,5).Value = "'86.88"
$newSheet.Cells.Item(3,5).Style = "Normal"
$newSheet.Cells.Item(3,6).Value = "'4.99"
$newSheet.Cells.Item(3,6).Style = "Normal"
$newSheet.Cells.Item(3,7).Value = "'0.0010"
$newSheet.Cells.Item(3,7).Style = "Normal"
$newSheet.Cells.Item(3,8).Value = 9

# ---------------------------------------------------------------------------
# 2) Prepend "2022-Q1" summary row to "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

# Restore column-A styling for the newly-inserted row (Insert() only carries
# formatting for columns that already had a cell in the row above).
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(2,1).PasteSpecial(-4122)
# Clear the copied formatting from the data cells on the new row.
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 2
$totalSheet.Cells.Item(2,4).Value = 0

# Renumber the index column for the rows pushed down.
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
